# luban excel export fixed
# The "Tags" column (M) on the BuffConfig sheet stored list values wrapped
# in square brackets (e.g. "[1]", "[2,3]"). Luban's exporter expects the
# brackets to be stripped - single-element lists collapse down to a bare
# integer, multi-element lists become a plain comma separated string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Single-value tags -> numeric
$ws.Range("M5").Value  = 1
$ws.Range("M6").Value  = 1
$ws.Range("M7").Value  = 1
$ws.Range("M8").Value  = 2
$ws.Range("M9").Value  = 2
$ws.Range("M10").Value = 1
$ws.Range("M12").Value = 2
$ws.Range("M13").Value = 2

# Multi-value tags -> comma separated text, brackets removed
$ws.Range("M11").Value = "2,3"
$ws.Range("M14").Value = "1,4"
$ws.Range("M15").Value = "1,4"
$ws.Range("M16").Value = "1,4"
$ws.Range("M17").Value = "1,4"

# Selection follows the edited Tags column, matching the authored workbook
$ws.Range("M4:M17").Select()
